# Natmi following Dr Hou advice
# Recomputes the Bmp4-Rgmb LR-pair table after adding the "sCs" sending
# cluster to the analysis: existing ECs/FAPs block stats are refreshed and
# a brand-new sCs block (rows 12-16) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp4"
$ws.Range("C2").Value = "Rgmb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.675695
$ws.Range("H2").Value = 26.027085
$ws.Range("I2").Value = 0.5592117158070719
$ws.Range("J2").Value = 0.5592117158070719
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.012448
$ws.Range("N2").Value = 27.037344
$ws.Range("O2").Value = 0.2887266436017198
$ws.Range("P2").Value = 0.2887266436017197
$ws.Range("Q2").Value = 78.18925005135999
$ws.Range("R2").Value = 703.7032504622399
$ws.Range("S2").Value = 0.1614593217677346
$ws.Range("T2").Value = 0.1614593217677346

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp4"
$ws.Range("C3").Value = "Rgmb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.675695
$ws.Range("H3").Value = 26.027085
$ws.Range("I3").Value = 0.5592117158070719
$ws.Range("J3").Value = 0.5592117158070719
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.682916
$ws.Range("N3").Value = 29.048748
$ws.Range("O3").Value = 0.3102060435696706
$ws.Range("P3").Value = 0.3102060435696705
$ws.Range("Q3").Value = 84.00602592662
$ws.Range("R3").Value = 756.05423333958
$ws.Range("S3").Value = 0.1734708538783188
$ws.Range("T3").Value = 0.1734708538783187

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp4"
$ws.Range("C4").Value = "Rgmb"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.675695
$ws.Range("H4").Value = 26.027085
$ws.Range("I4").Value = 0.5592117158070719
$ws.Range("J4").Value = 0.5592117158070719
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.223096
$ws.Range("N4").Value = 9.669288
$ws.Range("O4").Value = 0.1032564837085472
$ws.Range("P4").Value = 0.1032564837085472
$ws.Range("Q4").Value = 27.96259785172
$ws.Range("R4").Value = 251.66338066548
$ws.Range("S4").Value = 0.05774223542286164
$ws.Range("T4").Value = 0.05774223542286162

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bmp4"
$ws.Range("C5").Value = "Rgmb"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.675695
$ws.Range("H5").Value = 26.027085
$ws.Range("I5").Value = 0.5592117158070719
$ws.Range("J5").Value = 0.5592117158070719
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.120218333333334
$ws.Range("N5").Value = 9.360655000000001
$ws.Range("O5").Value = 0.09996065072307608
$ws.Range("P5").Value = 0.09996065072307606
$ws.Range("Q5").Value = 27.07006259340833
$ws.Range("R5").Value = 243.630563340675
$ws.Range("S5").Value = 0.05589916700404279
$ws.Range("T5").Value = 0.05589916700404278

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Bmp4"
$ws.Range("C6").Value = "Rgmb"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.675695
$ws.Range("H6").Value = 26.027085
$ws.Range("I6").Value = 0.5592117158070719
$ws.Range("J6").Value = 0.5592117158070719
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.175787666666667
$ws.Range("N6").Value = 18.527363
$ws.Range("O6").Value = 0.1978501783969864
$ws.Range("P6").Value = 0.1978501783969864
$ws.Range("Q6").Value = 53.57925018076167
$ws.Range("R6").Value = 482.213251626855
$ws.Range("S6").Value = 0.110640137734114
$ws.Range("T6").Value = 0.110640137734114

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp4"
$ws.Range("C7").Value = "Rgmb"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.316050666666666
$ws.Range("H7").Value = 18.948152
$ws.Range("I7").Value = 0.407115456505913
$ws.Range("J7").Value = 0.407115456505913
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.012448
$ws.Range("N7").Value = 27.037344
$ws.Range("O7").Value = 0.2887266436017198
$ws.Range("P7").Value = 0.2887266436017197
$ws.Range("Q7").Value = 56.92307819869866
$ws.Range("R7").Value = 512.307703788288
$ws.Range("S7").Value = 0.1175450793153342
$ws.Range("T7").Value = 0.1175450793153341

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Bmp4"
$ws.Range("C8").Value = "Rgmb"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.316050666666666
$ws.Range("H8").Value = 18.948152
$ws.Range("I8").Value = 0.407115456505913
$ws.Range("J8").Value = 0.407115456505913
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.682916
$ws.Range("N8").Value = 29.048748
$ws.Range("O8").Value = 0.3102060435696706
$ws.Range("P8").Value = 0.3102060435696705
$ws.Range("Q8").Value = 61.15778805707733
$ws.Range("R8").Value = 550.4200925136961
$ws.Range("S8").Value = 0.1262896750387596
$ws.Range("T8").Value = 0.1262896750387595

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Bmp4"
$ws.Range("C9").Value = "Rgmb"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.316050666666666
$ws.Range("H9").Value = 18.948152
$ws.Range("I9").Value = 0.407115456505913
$ws.Range("J9").Value = 0.407115456505913
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.223096
$ws.Range("N9").Value = 9.669288
$ws.Range("O9").Value = 0.1032564837085472
$ws.Range("P9").Value = 0.1032564837085472
$ws.Range("Q9").Value = 20.35723763953067
$ws.Range("R9").Value = 183.215138755776
$ws.Range("S9").Value = 0.04203731050220055
$ws.Range("T9").Value = 0.04203731050220054

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Bmp4"
$ws.Range("C10").Value = "Rgmb"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.316050666666666
$ws.Range("H10").Value = 18.948152
$ws.Range("I10").Value = 0.407115456505913
$ws.Range("J10").Value = 0.407115456505913
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.120218333333334
$ws.Range("N10").Value = 9.360655000000001
$ws.Range("O10").Value = 0.09996065072307608
$ws.Range("P10").Value = 0.09996065072307606
$ws.Range("Q10").Value = 19.70745708439556
$ws.Range("R10").Value = 177.36711375956
$ws.Range("S10").Value = 0.04069552595175324
$ws.Range("T10").Value = 0.04069552595175323

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Bmp4"
$ws.Range("C11").Value = "Rgmb"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6.316050666666666
$ws.Range("H11").Value = 18.948152
$ws.Range("I11").Value = 0.407115456505913
$ws.Range("J11").Value = 0.407115456505913
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.175787666666667
$ws.Range("N11").Value = 18.527363
$ws.Range("O11").Value = 0.1978501783969864
$ws.Range("P11").Value = 0.1978501783969864
$ws.Range("Q11").Value = 39.00658780924178
$ws.Range("R11").Value = 351.059290283176
$ws.Range("S11").Value = 0.08054786569786546
$ws.Range("T11").Value = 0.08054786569786543

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Bmp4"
$ws.Range("C12").Value = "Rgmb"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.5224053333333334
$ws.Range("H12").Value = 1.567216
$ws.Range("I12").Value = 0.03367282768701513
$ws.Range("J12").Value = 0.03367282768701512
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 9.012448
$ws.Range("N12").Value = 27.037344
$ws.Range("O12").Value = 0.2887266436017198
$ws.Range("P12").Value = 0.2887266436017197
$ws.Range("Q12").Value = 4.708150901589334
$ws.Range("R12").Value = 42.373358114304
$ws.Range("S12").Value = 0.009722242518650939
$ws.Range("T12").Value = 0.009722242518650934

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Bmp4"
$ws.Range("C13").Value = "Rgmb"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.5224053333333334
$ws.Range("H13").Value = 1.567216
$ws.Range("I13").Value = 0.03367282768701513
$ws.Range("J13").Value = 0.03367282768701512
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.682916
$ws.Range("N13").Value = 29.048748
$ws.Range("O13").Value = 0.3102060435696706
$ws.Range("P13").Value = 0.3102060435696705
$ws.Range("Q13").Value = 5.058406960618668
$ws.Range("R13").Value = 45.52566264556801
$ws.Range("S13").Value = 0.01044551465259222
$ws.Range("T13").Value = 0.01044551465259222

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Bmp4"
$ws.Range("C14").Value = "Rgmb"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.5224053333333334
$ws.Range("H14").Value = 1.567216
$ws.Range("I14").Value = 0.03367282768701513
$ws.Range("J14").Value = 0.03367282768701512
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.223096
$ws.Range("N14").Value = 9.669288
$ws.Range("O14").Value = 0.1032564837085472
$ws.Range("P14").Value = 0.1032564837085472
$ws.Range("Q14").Value = 1.683762540245334
$ws.Range("R14").Value = 15.153862862208
$ws.Range("S14").Value = 0.003476937783484994
$ws.Range("T14").Value = 0.003476937783484992

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Bmp4"
$ws.Range("C15").Value = "Rgmb"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.5224053333333334
$ws.Range("H15").Value = 1.567216
$ws.Range("I15").Value = 0.03367282768701513
$ws.Range("J15").Value = 0.03367282768701512
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.120218333333334
$ws.Range("N15").Value = 9.360655000000001
$ws.Range("O15").Value = 0.09996065072307608
$ws.Range("P15").Value = 0.09996065072307606
$ws.Range("Q15").Value = 1.630018698497778
$ws.Range("R15").Value = 14.67016828648
$ws.Range("S15").Value = 0.003365957767280045
$ws.Range("T15").Value = 0.003365957767280043

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Bmp4"
$ws.Range("C16").Value = "Rgmb"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.5224053333333334
$ws.Range("H16").Value = 1.567216
$ws.Range("I16").Value = 0.03367282768701513
$ws.Range("J16").Value = 0.03367282768701512
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 6.175787666666667
$ws.Range("N16").Value = 18.527363
$ws.Range("O16").Value = 0.1978501783969864
$ws.Range("P16").Value = 0.1978501783969864
$ws.Range("Q16").Value = 3.22626441460089
$ws.Range("R16").Value = 29.036379731408
$ws.Range("S16").Value = 0.006662174965006927
$ws.Range("T16").Value = 0.006662174965006924
